$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 17 (data now spans rows 2-16 instead of 2-17)
$ws.Rows.Item(17).Delete()

# Overwrite data rows 2-16 with the recomputed TPM values and updated cluster assignments
# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Fgf2"
$ws.Cells.Item(2, 3).Value = "Fgfr2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 1.066124666666667
$ws.Cells.Item(2, 8).Value = 3.198374
$ws.Cells.Item(2, 9).Value = 0.1044113535211941
$ws.Cells.Item(2, 10).Value = 0.1044113535211941
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.8155003333333334
$ws.Cells.Item(2, 14).Value = 2.446501
$ws.Cells.Item(2, 15).Value = 0.1910612426590028
$ws.Cells.Item(2, 16).Value = 0.1910612426590029
$ws.Cells.Item(2, 17).Value = 0.8694250210415557
$ws.Cells.Item(2, 18).Value = 7.824825189374001
$ws.Cells.Item(2, 19).Value = 0.0199489629514678
$ws.Cells.Item(2, 20).Value = 0.0199489629514678

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Fgf2"
$ws.Cells.Item(3, 3).Value = "Fgfr2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 1.066124666666667
$ws.Cells.Item(3, 8).Value = 3.198374
$ws.Cells.Item(3, 9).Value = 0.1044113535211941
$ws.Cells.Item(3, 10).Value = 0.1044113535211941
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 3.333134333333334
$ws.Cells.Item(3, 14).Value = 9.999403000000001
$ws.Cells.Item(3, 15).Value = 0.7809105179307759
$ws.Cells.Item(3, 16).Value = 0.780910517930776
$ws.Cells.Item(3, 17).Value = 3.553536730080223
$ws.Cells.Item(3, 18).Value = 31.98183057072201
$ws.Cells.Item(3, 19).Value = 0.08153592415608904
$ws.Cells.Item(3, 20).Value = 0.08153592415608905

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Fgf2"
$ws.Cells.Item(4, 3).Value = "Fgfr2"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 1.066124666666667
$ws.Cells.Item(4, 8).Value = 3.198374
$ws.Cells.Item(4, 9).Value = 0.1044113535211941
$ws.Cells.Item(4, 10).Value = 0.1044113535211941
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.119632
$ws.Cells.Item(4, 14).Value = 0.358896
$ws.Cells.Item(4, 15).Value = 0.02802823941022116
$ws.Cells.Item(4, 16).Value = 0.02802823941022117
$ws.Cells.Item(4, 17).Value = 0.1275426261226667
$ws.Cells.Item(4, 18).Value = 1.147883635104
$ws.Cells.Item(4, 19).Value = 0.002926466413637267
$ws.Cells.Item(4, 20).Value = 0.002926466413637267

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Fgf2"
$ws.Cells.Item(5, 3).Value = "Fgfr2"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 7.418580000000001
$ws.Cells.Item(5, 8).Value = 22.25574
$ws.Cells.Item(5, 9).Value = 0.7265416542955204
$ws.Cells.Item(5, 10).Value = 0.7265416542955204
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.8155003333333334
$ws.Cells.Item(5, 14).Value = 2.446501
$ws.Cells.Item(5, 15).Value = 0.1910612426590028
$ws.Cells.Item(5, 16).Value = 0.1910612426590029
$ws.Cells.Item(5, 17).Value = 6.049854462860002
$ws.Cells.Item(5, 18).Value = 54.44869016574001
$ws.Cells.Item(5, 19).Value = 0.1388139513132298
$ws.Cells.Item(5, 20).Value = 0.1388139513132298

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Fgf2"
$ws.Cells.Item(6, 3).Value = "Fgfr2"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 7.418580000000001
$ws.Cells.Item(6, 8).Value = 22.25574
$ws.Cells.Item(6, 9).Value = 0.7265416542955204
$ws.Cells.Item(6, 10).Value = 0.7265416542955204
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 3.333134333333334
$ws.Cells.Item(6, 14).Value = 9.999403000000001
$ws.Cells.Item(6, 15).Value = 0.7809105179307759
$ws.Cells.Item(6, 16).Value = 0.780910517930776
$ws.Cells.Item(6, 17).Value = 24.72712370258001
$ws.Cells.Item(6, 18).Value = 222.54411332322
$ws.Cells.Item(6, 19).Value = 0.5673640195541976
$ws.Cells.Item(6, 20).Value = 0.5673640195541977

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Fgf2"
$ws.Cells.Item(7, 3).Value = "Fgfr2"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 7.418580000000001
$ws.Cells.Item(7, 8).Value = 22.25574
$ws.Cells.Item(7, 9).Value = 0.7265416542955204
$ws.Cells.Item(7, 10).Value = 0.7265416542955204
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.119632
$ws.Cells.Item(7, 14).Value = 0.358896
$ws.Cells.Item(7, 15).Value = 0.02802823941022116
$ws.Cells.Item(7, 16).Value = 0.02802823941022117
$ws.Cells.Item(7, 17).Value = 0.8874995625600002
$ws.Cells.Item(7, 18).Value = 7.987496063040001
$ws.Cells.Item(7, 19).Value = 0.02036368342809299
$ws.Cells.Item(7, 20).Value = 0.02036368342809299

# Row 8
$ws.Cells.Item(8, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(8, 2).Value = "Fgf2"
$ws.Cells.Item(8, 3).Value = "Fgfr2"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.4336433333333334
$ws.Cells.Item(8, 8).Value = 1.30093
$ws.Cells.Item(8, 9).Value = 0.04246903649677213
$ws.Cells.Item(8, 10).Value = 0.04246903649677213
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.8155003333333334
$ws.Cells.Item(8, 14).Value = 2.446501
$ws.Cells.Item(8, 15).Value = 0.1910612426590028
$ws.Cells.Item(8, 16).Value = 0.1910612426590029
$ws.Cells.Item(8, 17).Value = 0.3536362828811112
$ws.Cells.Item(8, 18).Value = 3.18272654593
$ws.Cells.Item(8, 19).Value = 0.008114186887603828
$ws.Cells.Item(8, 20).Value = 0.00811418688760383

# Row 9
$ws.Cells.Item(9, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(9, 2).Value = "Fgf2"
$ws.Cells.Item(9, 3).Value = "Fgfr2"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.4336433333333334
$ws.Cells.Item(9, 8).Value = 1.30093
$ws.Cells.Item(9, 9).Value = 0.04246903649677213
$ws.Cells.Item(9, 10).Value = 0.04246903649677213
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 3.333134333333334
$ws.Cells.Item(9, 14).Value = 9.999403000000001
$ws.Cells.Item(9, 15).Value = 0.7809105179307759
$ws.Cells.Item(9, 16).Value = 0.780910517930776
$ws.Cells.Item(9, 17).Value = 1.445391482754445
$ws.Cells.Item(9, 18).Value = 13.00852334479
$ws.Cells.Item(9, 19).Value = 0.03316451728671535
$ws.Cells.Item(9, 20).Value = 0.03316451728671536

# Row 10
$ws.Cells.Item(10, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(10, 2).Value = "Fgf2"
$ws.Cells.Item(10, 3).Value = "Fgfr2"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.4336433333333334
$ws.Cells.Item(10, 8).Value = 1.30093
$ws.Cells.Item(10, 9).Value = 0.04246903649677213
$ws.Cells.Item(10, 10).Value = 0.04246903649677213
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.119632
$ws.Cells.Item(10, 14).Value = 0.358896
$ws.Cells.Item(10, 15).Value = 0.02802823941022116
$ws.Cells.Item(10, 16).Value = 0.02802823941022117
$ws.Cells.Item(10, 17).Value = 0.05187761925333334
$ws.Cells.Item(10, 18).Value = 0.4668985732800001
$ws.Cells.Item(10, 19).Value = 0.00119033232245295
$ws.Cells.Item(10, 20).Value = 0.00119033232245295

# Row 11
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Fgf2"
$ws.Cells.Item(11, 3).Value = "Fgfr2"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.951285
$ws.Cells.Item(11, 8).Value = 2.853855
$ws.Cells.Item(11, 9).Value = 0.09316448398568379
$ws.Cells.Item(11, 10).Value = 0.09316448398568379
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0.8155003333333334
$ws.Cells.Item(11, 14).Value = 2.446501
$ws.Cells.Item(11, 15).Value = 0.1910612426590028
$ws.Cells.Item(11, 16).Value = 0.1910612426590029
$ws.Cells.Item(11, 17).Value = 0.7757732345950001
$ws.Cells.Item(11, 18).Value = 6.981959111355001
$ws.Cells.Item(11, 19).Value = 0.01780012208198951
$ws.Cells.Item(11, 20).Value = 0.01780012208198952

# Row 12
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Fgf2"
$ws.Cells.Item(12, 3).Value = "Fgfr2"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.951285
$ws.Cells.Item(12, 8).Value = 2.853855
$ws.Cells.Item(12, 9).Value = 0.09316448398568379
$ws.Cells.Item(12, 10).Value = 0.09316448398568379
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 3.333134333333334
$ws.Cells.Item(12, 14).Value = 9.999403000000001
$ws.Cells.Item(12, 15).Value = 0.7809105179307759
$ws.Cells.Item(12, 16).Value = 0.780910517930776
$ws.Cells.Item(12, 17).Value = 3.170760694285001
$ws.Cells.Item(12, 18).Value = 28.53684624856501
$ws.Cells.Item(12, 19).Value = 0.0727531254420138
$ws.Cells.Item(12, 20).Value = 0.07275312544201382

# Row 13
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Fgf2"
$ws.Cells.Item(13, 3).Value = "Fgfr2"
$ws.Cells.Item(13, 4).Value = "MuSCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.951285
$ws.Cells.Item(13, 8).Value = 2.853855
$ws.Cells.Item(13, 9).Value = 0.09316448398568379
$ws.Cells.Item(13, 10).Value = 0.09316448398568379
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.119632
$ws.Cells.Item(13, 14).Value = 0.358896
$ws.Cells.Item(13, 15).Value = 0.02802823941022116
$ws.Cells.Item(13, 16).Value = 0.02802823941022117
$ws.Cells.Item(13, 17).Value = 0.11380412712
$ws.Cells.Item(13, 18).Value = 1.02423714408
$ws.Cells.Item(13, 19).Value = 0.002611236461680461
$ws.Cells.Item(13, 20).Value = 0.002611236461680461

# Row 14
$ws.Cells.Item(14, 1).Value = "Resolving-Mac"
$ws.Cells.Item(14, 2).Value = "Fgf2"
$ws.Cells.Item(14, 3).Value = "Fgfr2"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.3411786666666667
$ws.Cells.Item(14, 8).Value = 1.023536
$ws.Cells.Item(14, 9).Value = 0.03341347170082953
$ws.Cells.Item(14, 10).Value = 0.03341347170082953
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.8155003333333334
$ws.Cells.Item(14, 14).Value = 2.446501
$ws.Cells.Item(14, 15).Value = 0.1910612426590028
$ws.Cells.Item(14, 16).Value = 0.1910612426590029
$ws.Cells.Item(14, 17).Value = 0.2782313163928889
$ws.Cells.Item(14, 18).Value = 2.504081847536
$ws.Cells.Item(14, 19).Value = 0.006384019424711914
$ws.Cells.Item(14, 20).Value = 0.006384019424711915

# Row 15
$ws.Cells.Item(15, 1).Value = "Resolving-Mac"
$ws.Cells.Item(15, 2).Value = "Fgf2"
$ws.Cells.Item(15, 3).Value = "Fgfr2"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.3411786666666667
$ws.Cells.Item(15, 8).Value = 1.023536
$ws.Cells.Item(15, 9).Value = 0.03341347170082953
$ws.Cells.Item(15, 10).Value = 0.03341347170082953
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 3.333134333333334
$ws.Cells.Item(15, 14).Value = 9.999403000000001
$ws.Cells.Item(15, 15).Value = 0.7809105179307759
$ws.Cells.Item(15, 16).Value = 0.780910517930776
$ws.Cells.Item(15, 17).Value = 1.137194327667556
$ws.Cells.Item(15, 18).Value = 10.234748949008
$ws.Cells.Item(15, 19).Value = 0.02609293149176011
$ws.Cells.Item(15, 20).Value = 0.02609293149176012

# Row 16
$ws.Cells.Item(16, 1).Value = "Resolving-Mac"
$ws.Cells.Item(16, 2).Value = "Fgf2"
$ws.Cells.Item(16, 3).Value = "Fgfr2"
$ws.Cells.Item(16, 4).Value = "MuSCs"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.3411786666666667
$ws.Cells.Item(16, 8).Value = 1.023536
$ws.Cells.Item(16, 9).Value = 0.03341347170082953
$ws.Cells.Item(16, 10).Value = 0.03341347170082953
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 0.119632
$ws.Cells.Item(16, 14).Value = 0.358896
$ws.Cells.Item(16, 15).Value = 0.02802823941022116
$ws.Cells.Item(16, 16).Value = 0.02802823941022117
$ws.Cells.Item(16, 17).Value = 0.04081588625066667
$ws.Cells.Item(16, 18).Value = 0.367342976256
$ws.Cells.Item(16, 19).Value = 0.0009365207843574997
$ws.Cells.Item(16, 20).Value = 0.0009365207843574998
